$d = $word.ActiveDocument

# Hex color 2C3E50 expressed as the BGR integer Word's Font.Color expects
# (0x2C, 0x3E, 0x50) -> bb*65536 + gg*256 + rr
$HighlightColor = 0x50 * 65536 + 0x3E * 256 + 0x2C

# Finds $searchText inside the range [$fromPos, end-of-paragraph], makes it
# bold + colored, and returns the position right after the match so the next
# search in the same paragraph starts from there (handles repeated
# substrings like " to " safely).
function Format-Segment($paraRange, $searchText, $fromPos) {
    $pEnd = $paraRange.End
    $rng = $d.Range($fromPos, $pEnd)
    $found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = 1
        $rng.Font.Color = $HighlightColor
        return $rng.End
    }
    Write-Output "NOT FOUND: $searchText"
    return $fromPos
}

# Returns the first paragraph whose text matches the -like wildcard pattern.
function Get-ParaByText($marker) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like $marker) {
            return $p
        }
    }
    return $null
}

# 1) "...developed geospatial machine learning algorithms improving demographic
#     classification accuracy from 23% to 64%"
$p = Get-ParaByText "*Discovered systematic race coding errors affecting all Black and Asian-American voters, developed*"
if ($p -ne $null) {
    $pr = $p.Range
    $pos = $pr.Start
    $pos = Format-Segment $pr "23%" $pos
    $pos = Format-Segment $pr "64%" $pos
}

# 2) "Achieved 87% prediction accuracy for voter turnout vs. industry standard
#     of 71%, reducing polling error margins from ±4.2% to ±2.1%"
$p = Get-ParaByText "*Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing*"
if ($p -ne $null) {
    $pr = $p.Range
    $pos = $pr.Start
    $pos = Format-Segment $pr "87%" $pos
    $pos = Format-Segment $pr "71%" $pos
    $pos = Format-Segment $pr "±4.2%" $pos
    $pos = Format-Segment $pr "±2.1%" $pos
}

# 3) "Wrote RFP and analyzed bids from 1,200 vendors for research platform
#     development"
$p = Get-ParaByText "*Wrote RFP and analyzed bids from 1,200 vendors*"
if ($p -ne $null) {
    $pr = $p.Range
    $pos = $pr.Start
    $pos = Format-Segment $pr "1,200" $pos
}

# 4) "Created comprehensive meta-analysis framework handling millions of survey
#     responses that became the $400M Polling Consortium Database at The
#     Analyst Institute, now valued at $1B+"
$p = Get-ParaByText "*Created comprehensive meta-analysis framework handling millions of survey responses*"
if ($p -ne $null) {
    $pr = $p.Range
    $pos = $pr.Start
    $pos = Format-Segment $pr "`$400M" $pos
    $pos = Format-Segment $pr "`$1B" $pos
}

# 5) "Algorithm reduced mapping costs by 73.5%, saving campaigns and
#     organizations $4.7M"
$p = Get-ParaByText "*Algorithm reduced mapping costs by 73.5%*"
if ($p -ne $null) {
    $pr = $p.Range
    $pos = $pr.Start
    $pos = Format-Segment $pr "73.5%" $pos
    $pos = Format-Segment $pr "`$4.7M" $pos
}

# 6) "Achieved 87% prediction accuracy for voter turnout vs. industry standard
#     of 71%" (short form, without the polling-error-margin clause). Paragraph
#     text always ends with a trailing paragraph-mark character, so match
#     without anchoring to the end of the wildcard and exclude the longer
#     "...reducing polling error margins..." variant explicitly.
$p = $null
foreach ($cand in $d.Paragraphs) {
    $ct = $cand.Range.Text
    if ($ct -like "*Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%*" -and $ct -notlike "*reducing*") {
        $p = $cand
        break
    }
}
if ($p -ne $null) {
    $pr = $p.Range
    $pos = $pr.Start
    $pos = Format-Segment $pr "87%" $pos
    $pos = Format-Segment $pr "71%" $pos
}

Write-Output "done"
